$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows (row1 = headers: firstname, middlename, lastname - unchanged)
# Row 2: firstname/middlename/lastname
$ws.Range("A2").Value = "Jay"
$ws.Range("B2").Value = "Peter"
$ws.Range("C2").Value = "Jackson"

# Row 3
$ws.Range("A3").Value = "Alison"
$ws.Range("B3").Value = "Marie"
$ws.Range("C3").Value = "Palm"

# Row 4
$ws.Range("A4").Value = "Connor"
$ws.Range("B4").Value = "Davia"
$ws.Range("C4").Value = "Jones"
